# "changes in the reset page"
#  - Adds a new worksheet "ResetPage" after "LoginPage"
#  - Populates it with an "expectedmessage" / error-message row, formatted
#    like the matching expected-value cell on LoginPage
#  - Makes ResetPage the active/selected sheet (tabSelected + activeTab move
#    from LoginPage to ResetPage)

$wb = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("LoginPage")

# Add the new worksheet at the end of the workbook and name it.
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsNew.Name = "ResetPage"

# Fill in the data for the reset page.
$wsNew.Range("A1").Value = "expectedmessage"
$wsNew.Range("B1").Value = "We can't find a user with that e-mail address."

# Match the formatting used for the equivalent "expected value" cell on
# LoginPage (B1, the small Consolas expected-value style) so the reset
# message cell looks the same as the other expected-value cells.
$wsLogin.Range("B1").Copy()
$wsNew.Range("B1").PasteSpecial(-4122)

# Column widths for the new sheet.
$wsNew.Columns.Item(1).ColumnWidth = 19.5
$wsNew.Columns.Item(2).ColumnWidth = 40

# Select a cell on the new sheet (matches the stored selection).
$wsNew.Range("C6").Select()

# Make the new sheet the active tab - moves tabSelected/activeTab here and
# off of LoginPage.
$wsNew.Activate()
